# Update automatico via Actualizar 03-09-2021 12-34-14
#
# Column D ("Fecha") holds a rolling "last updated" timestamp shared by
# blocks of 14 rows. On each automatic refresh the newest timestamp is
# pushed into the top block and every older block's value shifts down
# to the next block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTop = 44264.52353865721
$newMid = 44264.50214331019
$newBot = 44264.4807412963

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newTop
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $newMid
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $newBot
}
